$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)
$ws4 = $wb.Worksheets.Item(4)

# Sheet 1
$ws1.Range("F5").Value = 326
$ws1.Range("F6").Value = 470
$ws1.Range("F7").Value = 2154
$ws1.Range("F9").Value = 52
$ws1.Range("F10").Value = 1630
$ws1.Range("F11").Value = 1630
$ws1.Range("F12").Value = 1360
$ws1.Range("F17").Value = 576
$ws1.Range("F18").Value = 157
$ws1.Range("F19").Value = 11
$ws1.Range("F20").Value = 7231
$ws1.Range("F21").Value = 7930
$ws1.Range("F24").Value = 196
$ws1.Range("F33").Value = 275
$ws1.Range("F35").Value = 1436
$ws1.Range("F36").Value = 206
$ws1.Range("F39").Value = 291
$ws1.Range("F41").Value = 724
$ws1.Range("F43").Value = 1363
$ws1.Range("F44").Value = 341
$ws1.Range("F46").Value = 193
$ws1.Range("F48").Value = 175

# Sheet 2
$ws2.Range("F9").Value = 26

# Sheet 3
$ws3.Range("F3").Value = 2621
$ws3.Range("F4").Value = 281
$ws3.Range("F6").Value = 12

# Sheet 4
$ws4.Range("F7").Value = 326
$ws4.Range("F9").Value = 470
$ws4.Range("F10").Value = 2154
$ws4.Range("F12").Value = 52
$ws4.Range("F13").Value = 1630
$ws4.Range("F14").Value = 1630
$ws4.Range("F18").Value = 576
$ws4.Range("F20").Value = 157
$ws4.Range("F23").Value = 11
$ws4.Range("F24").Value = 7231
$ws4.Range("F25").Value = 7930
$ws4.Range("F30").Value = 1436
$ws4.Range("F31").Value = 206
$ws4.Range("F35").Value = 291
$ws4.Range("F37").Value = 26
$ws4.Range("F39").Value = 724
$ws4.Range("F43").Value = 1363
$ws4.Range("F44").Value = 341
$ws4.Range("F46").Value = 193
$ws4.Range("F47").Value = 175
